# Beneficiarios workbook update: remove the last beneficiary row (CURP
# AAAE560802MZSLVR07 / ERNESTINA / C DE LA FE) together with the trailing
# helper/formula row, renumber the "Clave Municipio" column, and fix up the
# hyperlink range + active selection that shift as a result.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the 10th beneficiary row (row 11) and the formula helper row
# (row 12) that followed it. Deleting both at once lets Excel shift the
# remaining rows/formulas up in a single operation.
$ws.Range("A11:AM12").EntireRow.Delete()

# The hyperlink range that covered AH3:AH11 needs to shrink to AH3:AH10
# now that the last data row is row 10. Recreate it pointing at the same
# mailto target used by the existing hyperlinks.
$ws.Hyperlinks.Item(2).Delete()
$ws.Hyperlinks.Add($ws.Range("AH3:AH10"), "mailto:example@hotmail.com", "", "", "example@hotmail.com")

# Fill in the "Clave Municipio" column (AM) for the 9 remaining
# beneficiary rows with a simple 1-based sequence.
for ($i = 0; $i -le 8; $i++) {
    $row = 2 + $i
    $ws.Range("AM$row").Value = $i + 1
}

# Update the active selection to match the post-edit view.
$ws.Range("AN9").Select()
